$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.479.19"
$ws.Range("E2").Value = "  -1.67%  "

$ws.Range("D3").Value = "2.631.71"
$ws.Range("E3").Value = "  -1.30%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.98%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.15%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.651"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.95%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  -3.68%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.81"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.05%  "

$ws.Range("E11").Value = "  -2.03%  "

$ws.Range("E12").Value = "  +0.16%  "

$ws.Range("E13").Value = "  -1.25%  "

$ws.Range("E14").Value = "  -5.46%  "

$ws.Range("D15").Value = "3.104.30"
$ws.Range("E15").Value = "  -1.34%  "

$ws.Range("D16").Value = "64.283.80"
$ws.Range("E16").Value = "  -1.86%  "

$ws.Range("D17").Value = "2.630.62"
$ws.Range("E17").Value = "  -1.17%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.27"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.91%  "

$ws.Range("E19").Value = "  -2.37%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.85%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "346.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.15%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.15%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.54%  "

$ws.Range("E24").Value = "  +0.55%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.75"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.99%  "

$ws.Range("E26").Value = "  -3.00%  "

$ws.Range("E27").Value = "  -1.26%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "567.59"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.02%  "

$ws.Range("E29").Value = "  -1.82%  "

$ws.Range("E30").Value = "  -0.01%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.61%  "

$ws.Range("E32").Value = "  -2.22%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.68"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.72%  "

$ws.Range("E34").Value = "  -3.45%  "

$ws.Range("E35").Value = "  -2.36%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.413"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.93%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.06"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.45%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.998"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.14%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.93"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.39%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "154.79"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.31%  "

$ws.Range("E41").Value = "  -0.01%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.47"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.29%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "158.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.49%  "

$ws.Range("E44").Value = "  -2.25%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "23.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.97%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0599"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.43%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.636"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.66%  "

$ws.Range("E48").Value = "  +4.05%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0252"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.21%  "

$ws.Range("E50").Value = "  -3.24%  "

$ws.Range("D51").Value = "0.0₆0238"
$ws.Range("E51").Value = "  -5.75%  "
